# Insert a new row at position 198, pushing existing rows 198:249 down to 199:250.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(198).Insert()

# Populate the newly inserted row 198 with the new record.
$ws.Range("A198").Value = 4
$ws.Range("B198").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C198").Value = "Los Lagos"
$ws.Range("D198").Value = 44642
$ws.Range("E198").Value = 10
$ws.Range("F198").Value = 100112037
$ws.Range("G198").Value = "Cebollín"
$ws.Range("H198").Value = "Sin especificar"
$ws.Range("I198").Value = "Segunda"
$ws.Range("J198").Value = 160
$ws.Range("K198").Value = 7500
$ws.Range("L198").Value = 8500
$ws.Range("M198").Value = 8000
$ws.Range("N198").Value = "$/paquete 36 unidades"
$ws.Range("O198").Value = "Región Metropolitana"
$ws.Range("P198").Value = 222
$ws.Range("Q198").Value = 36
$ws.Range("R198").Value = "Hortaliza"
